# Add 2022-Q3 data
#
# 1. Insert a new worksheet "2022-Q3" right after the "总计" (summary) sheet,
#    holding the fund-holding detail table for the new quarter.
# 2. Insert a new row into the "总计" summary sheet for "2022-Q3", pushing
#    the existing quarterly rows down by one and renumbering the index
#    column.

$wb = $excel.ActiveWorkbook
$total = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------
# Step 1: new "2022-Q3" worksheet (fund holdings detail)
# ---------------------------------------------------------------------
$newSheet = $wb.Worksheets.Add($null, $total)
$newSheet.Name = "2022-Q3"

$newSheet.Range("B1").Value2 = "基金代码"
$newSheet.Range("C1").Value2 = "基金名称"
$newSheet.Range("D1").Value2 = "基金规模"
$newSheet.Range("E1").Value2 = "股票总仓位"
$newSheet.Range("F1").Value2 = "仓位占比"
$newSheet.Range("G1").Value2 = "持有市值(亿元)"
$newSheet.Range("H1").Value2 = "仓位排名"

# match the bold/centered header style used by every other sheet
$total.Range("B1:D1").Copy()
$newSheet.Range("B1:H1").PasteSpecial(-4122, $null, $false, $false)

$newSheet.Range("A2").Value2 = 0
$total.Range("A2").Copy()
$newSheet.Range("A2").PasteSpecial(-4122, $null, $false, $false)

$newSheet.Range("B2").Value2 = "'003413"
$newSheet.Range("B2").Style = "Normal"
$newSheet.Range("C2").Value2 = "华泰柏瑞新经济沪港深混合"
$newSheet.Range("D2").Value2 = "'0.42"
$newSheet.Range("D2").Style = "Normal"
$newSheet.Range("E2").Value2 = "'86.45"
$newSheet.Range("E2").Style = "Normal"
$newSheet.Range("F2").Value2 = "'5.44"
$newSheet.Range("F2").Style = "Normal"
$newSheet.Range("G2").Value2 = "'0.0228"
$newSheet.Range("G2").Style = "Normal"
$newSheet.Range("H2").Value2 = 7

# ---------------------------------------------------------------------
# Step 2: push rows 2-6 of "总计" down to 3-7, then fill in the new row 2
# ---------------------------------------------------------------------
for ($r = 6; $r -ge 2; $r--) {
    $dst = $r + 1
    $total.Cells.Item($dst, 2).Value2 = $total.Cells.Item($r, 2).Value2
    $total.Cells.Item($dst, 3).Value2 = $total.Cells.Item($r, 3).Value2
    $total.Cells.Item($dst, 4).Value2 = $total.Cells.Item($r, 4).Value2
}

# row 7 is brand-new on this sheet - give its index cell the same style as
# the rest of column A
$total.Range("A6").Copy()
$total.Range("A7").PasteSpecial(-4122, $null, $false, $false)

$total.Range("B2").Value2 = "2022-Q3"
$total.Range("C2").Value2 = 1
$total.Range("D2").Value2 = 0.02

# renumber the running index in column A (0..5)
for ($r = 2; $r -le 7; $r++) {
    $total.Cells.Item($r, 1).Value2 = $r - 2
}
